# Apply the crypto-price-table refresh described by the commit diff.
# Numeric-looking text values (e.g. "206.36") are written with a leading
# apostrophe so Excel stores them as TEXT (matching the workbook's original
# inlineStr cells) instead of auto-converting them to numbers; ClearFormats()
# then drops the quote-prefix style Excel applies, so no stray cell style is
# introduced and the cell keeps its original (default) formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.927.54'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.551.58'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").Value = '''206.36'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = '''0.488'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").Value = '''21.94'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '''0.0593'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").Value = '1.772.78'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = '1.554.78'
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '26.907.53'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '0.0₃0711'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").Value = '''216.86'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '''7.29'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").Value = '''4.08'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("D25").Value = '''153.64'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").Value = '''14.96'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").Value = '''0.0468'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").Value = '''3.21'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").Value = '''3.11'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.39%  '
$ws.Range("D34").Value = '1.409.33'
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("D36").Value = '''0.962'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").Value = '''0.0166'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").Value = '''0.525'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  +3.08%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.30'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.30%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '''0.996'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").Value = '''64.45'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("D47").Value = '1.686.63'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").Value = '''87.34'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("E49").Value = '  +1.18%  '
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").Value = '''0.0957'
$ws.Range("D51").ClearFormats()
